# Weekly update of fruit/vegetable (Alcachofa) price records.
# The rows of data (2-16) get their values re-shuffled to reflect
# the latest weekly source data; only the cells whose values actually
# change are touched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, "D").Value = 44841
$ws.Cells.Item(2, "J").Value = 45
$ws.Cells.Item(2, "K").Value = 12000
$ws.Cells.Item(2, "L").Value = 12000
$ws.Cells.Item(2, "M").Value = 12000
$ws.Cells.Item(2, "P").Value = 400

# Row 3
$ws.Cells.Item(3, "D").Value = 44841
$ws.Cells.Item(3, "I").Value = "Segunda"
$ws.Cells.Item(3, "J").Value = 45
$ws.Cells.Item(3, "K").Value = 10000
$ws.Cells.Item(3, "L").Value = 10000
$ws.Cells.Item(3, "M").Value = 10000
$ws.Cells.Item(3, "N").Value = "$/caja 40 unidades"
$ws.Cells.Item(3, "O").Value = "Provincia de Limarí"
$ws.Cells.Item(3, "P").Value = 250
$ws.Cells.Item(3, "Q").Value = 40

# Row 4
$ws.Cells.Item(4, "D").Value = 44425
$ws.Cells.Item(4, "J").Value = 35

# Row 5
$ws.Cells.Item(5, "D").Value = 44435
$ws.Cells.Item(5, "J").Value = 25
$ws.Cells.Item(5, "K").Value = 14000
$ws.Cells.Item(5, "L").Value = 14000
$ws.Cells.Item(5, "M").Value = 14000
$ws.Cells.Item(5, "P").Value = 467

# Row 6
$ws.Cells.Item(6, "D").Value = 44435

# Row 7
$ws.Cells.Item(7, "D").Value = 44474
$ws.Cells.Item(7, "K").Value = 10000
$ws.Cells.Item(7, "L").Value = 10000
$ws.Cells.Item(7, "M").Value = 10000
$ws.Cells.Item(7, "P").Value = 333

# Row 8
$ws.Cells.Item(8, "D").Value = 44418
$ws.Cells.Item(8, "I").Value = "Primera"
$ws.Cells.Item(8, "J").Value = 30
$ws.Cells.Item(8, "K").Value = 15000
$ws.Cells.Item(8, "L").Value = 15000
$ws.Cells.Item(8, "M").Value = 15000
$ws.Cells.Item(8, "N").Value = "$/caja 30 unidades"
$ws.Cells.Item(8, "P").Value = 500
$ws.Cells.Item(8, "Q").Value = 30

# Row 9
$ws.Cells.Item(9, "D").Value = 44421
$ws.Cells.Item(9, "J").Value = 25
$ws.Cells.Item(9, "K").Value = 15000
$ws.Cells.Item(9, "L").Value = 16000
$ws.Cells.Item(9, "M").Value = 15400
$ws.Cells.Item(9, "P").Value = 513

# Row 10
$ws.Cells.Item(10, "D").Value = 44432
$ws.Cells.Item(10, "K").Value = 14000
$ws.Cells.Item(10, "L").Value = 14000
$ws.Cells.Item(10, "M").Value = 14000
$ws.Cells.Item(10, "O").Value = "Provincia del Elquí"
$ws.Cells.Item(10, "P").Value = 467

# Row 11
$ws.Cells.Item(11, "D").Value = 44460
$ws.Cells.Item(11, "J").Value = 45
$ws.Cells.Item(11, "K").Value = 13000
$ws.Cells.Item(11, "L").Value = 13000
$ws.Cells.Item(11, "M").Value = 13000
$ws.Cells.Item(11, "P").Value = 433

# Row 12
$ws.Cells.Item(12, "D").Value = 44467
$ws.Cells.Item(12, "J").Value = 35
$ws.Cells.Item(12, "K").Value = 12000
$ws.Cells.Item(12, "L").Value = 12000
$ws.Cells.Item(12, "M").Value = 12000
$ws.Cells.Item(12, "P").Value = 400

# Row 15
$ws.Cells.Item(15, "D").Value = 44449
$ws.Cells.Item(15, "J").Value = 45
$ws.Cells.Item(15, "K").Value = 12000
$ws.Cells.Item(15, "L").Value = 12000
$ws.Cells.Item(15, "M").Value = 12000
$ws.Cells.Item(15, "P").Value = 400

# Row 16
$ws.Cells.Item(16, "D").Value = 44446
$ws.Cells.Item(16, "J").Value = 25
$ws.Cells.Item(16, "K").Value = 14000
$ws.Cells.Item(16, "L").Value = 14000
$ws.Cells.Item(16, "M").Value = 14000
$ws.Cells.Item(16, "P").Value = 467
